$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "244.67"
Set-TextValue $ws "G2" "21"

# Row 3
Set-TextValue $ws "D3" "21.87"
Set-TextValue $ws "G3" "21"

# Row 4
Set-TextValue $ws "D4" "5.403"
Set-TextValue $ws "G4" "21"

# Row 5
Set-TextValue $ws "D5" "0.06038"
Set-TextValue $ws "G5" "21"

# Row 6
Set-TextValue $ws "D6" "3.394"
Set-TextValue $ws "G6" "21"

# Row 7
Set-TextValue $ws "D7" "0.8140"
Set-TextValue $ws "G7" "21"

# Row 8
Set-TextValue $ws "G8" "21"

# Row 9
Set-TextValue $ws "D9" "0.1434"
Set-TextValue $ws "G9" "21"

# Row 10
Set-TextValue $ws "D10" "0.07471"
Set-TextValue $ws "G10" "21"

# Row 11
Set-TextValue $ws "D11" "0.03387"
Set-TextValue $ws "G11" "21"

# Row 12
Set-TextValue $ws "D12" "0.03067"
Set-TextValue $ws "G12" "21"

# Row 13
Set-TextValue $ws "G13" "21"

# Row 14
Set-TextValue $ws "D14" "4.010"
Set-TextValue $ws "G14" "21"

# Row 15
Set-TextValue $ws "G15" "21"

# Row 16
Set-TextValue $ws "D16" "0.04814"
Set-TextValue $ws "G16" "21"

# Row 17
Set-TextValue $ws "G17" "21"

# Row 18
Set-TextValue $ws "D18" "0.005600"
Set-TextValue $ws "G18" "21"

# Row 19
Set-TextValue $ws "D19" "0.004162"
Set-TextValue $ws "G19" "21"

# Row 20
Set-TextValue $ws "D20" "0.0009872"
Set-TextValue $ws "G20" "21"

# Row 21
Set-TextValue $ws "D21" "3.665"
Set-TextValue $ws "G21" "21"

# Row 22
Set-TextValue $ws "D22" "6.428"
Set-TextValue $ws "G22" "21"

# Row 23
Set-TextValue $ws "G23" "21"

# Row 24
Set-TextValue $ws "G24" "21"

# Row 25
Set-TextValue $ws "G25" "21"

# Row 26
Set-TextValue $ws "G26" "21"

# Row 27
Set-TextValue $ws "D27" "0.0002901"
Set-TextValue $ws "G27" "21"

# Row 28
Set-TextValue $ws "G28" "21"

# Row 29
Set-TextValue $ws "G29" "21"

# Row 30
Set-TextValue $ws "G30" "21"

# Row 31
Set-TextValue $ws "G31" "21"

# Row 32
Set-TextValue $ws "G32" "21"

# Row 33
Set-TextValue $ws "G33" "21"

# Row 34
Set-TextValue $ws "G34" "21"

# Row 35
Set-TextValue $ws "G35" "21"

# Row 36
Set-TextValue $ws "G36" "21"

# Row 37
Set-TextValue $ws "G37" "21"

# Row 38
Set-TextValue $ws "G38" "21"

# Row 39
Set-TextValue $ws "G39" "21"

# Row 40
Set-TextValue $ws "D40" "0.03997"
Set-TextValue $ws "G40" "21"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006421"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue $ws "G41" "21"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1077"
$ws.Range("E42").Value = "41BKEXTokenBKK"
Set-TextValue $ws "G42" "21"

# Row 43
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002901"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws "G43" "21"

# Row 44
Set-TextValue $ws "D44" "0.005792"
Set-TextValue $ws "G44" "21"

# Row 45
Set-TextValue $ws "D45" "0.00005237"
Set-TextValue $ws "G45" "21"

# Row 46
Set-TextValue $ws "G46" "21"

# Row 47
Set-TextValue $ws "G47" "21"

# Row 48
Set-TextValue $ws "D48" "0.002322"
Set-TextValue $ws "G48" "21"

# Row 49
Set-TextValue $ws "G49" "21"

# Row 50
Set-TextValue $ws "G50" "21"

# Row 51
Set-TextValue $ws "G51" "21"
